$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 94 (Leve Item ID 19905)
$ws.Range("H94").Value = 883.9167
$ws.Range("I94").Value = 1002.1111
$ws.Range("K94").Value = 1002.1111
$ws.Range("M94").Value = -551.1111

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1499.381
$ws.Range("I98").Value = 1327.8889
$ws.Range("K98").Value = 1327.8889
$ws.Range("M98").Value = 170.1111000000001

# Row 101 (Leve Item ID 19884)
$ws.Range("H101").Value = 1050.2
$ws.Range("I101").Value = 812.75
$ws.Range("J101").Value = 2000
$ws.Range("K101").Value = 2438.25
$ws.Range("L101").Value = 6000
$ws.Range("M101").Value = -816.25
$ws.Range("N101").Value = -9244

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1499.381
$ws.Range("I122").Value = 1327.8889
$ws.Range("K122").Value = 3983.6667
$ws.Range("M122").Value = -1533.6667

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1305.625
$ws.Range("I132").Value = 1242.3611
$ws.Range("K132").Value = 3727.0833
$ws.Range("M132").Value = -1197.0833

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 9000.532999999999
$ws.Range("I137").Value = 4028.1956
$ws.Range("J137").Value = 16887.69
$ws.Range("K137").Value = 12084.5868
$ws.Range("L137").Value = 50663.06999999999
$ws.Range("M137").Value = -9534.586800000001
$ws.Range("N137").Value = -55763.06999999999

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 5151.623
$ws.Range("I138").Value = 4231.5
$ws.Range("J138").Value = 5947.4053
$ws.Range("K138").Value = 12694.5
$ws.Range("L138").Value = 17842.2159
$ws.Range("M138").Value = -7554.5
$ws.Range("N138").Value = -28122.2159

# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 2103.5557
$ws.Range("I141").Value = 2284.3333
$ws.Range("K141").Value = 6852.999899999999
$ws.Range("M141").Value = -1672.999899999999


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 5078.231
$ws.Range("I2").Value = 4496.8667
$ws.Range("J2").Value = 7016.1113
$ws.Range("K2").Value = 4496.8667
$ws.Range("L2").Value = 7016.1113
$ws.Range("M2").Value = -4383.8667
$ws.Range("N2").Value = -7242.1113

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 5959.231
$ws.Range("J74").Value = 6397.75
$ws.Range("L74").Value = 6397.75
$ws.Range("N74").Value = -8145.75

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 5959.231
$ws.Range("J77").Value = 6397.75
$ws.Range("L77").Value = 31988.75
$ws.Range("N77").Value = -40724.75

# Row 101 (Leve Item ID 18518)
$ws.Range("H101").Value = 47800
$ws.Range("J101").Value = 47800
$ws.Range("L101").Value = 47800
$ws.Range("N101").Value = -54290

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 1426.5238
$ws.Range("I102").Value = 1400.35
$ws.Range("K102").Value = 1400.35
$ws.Range("M102").Value = 221.6500000000001

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 5078.231
$ws.Range("I116").Value = 4496.8667
$ws.Range("J116").Value = 7016.1113
$ws.Range("K116").Value = 4496.8667
$ws.Range("L116").Value = 7016.1113
$ws.Range("M116").Value = -2202.8667
$ws.Range("N116").Value = -11604.1113

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 5686.1333
$ws.Range("I132").Value = 4087.818
$ws.Range("J132").Value = 10081.5
$ws.Range("K132").Value = 12263.454
$ws.Range("L132").Value = 30244.5
$ws.Range("M132").Value = -9733.454000000002
$ws.Range("N132").Value = -35304.5


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 5078.231
$ws.Range("I3").Value = 4496.8667
$ws.Range("J3").Value = 7016.1113
$ws.Range("K3").Value = 4496.8667
$ws.Range("L3").Value = 7016.1113
$ws.Range("M3").Value = -4382.8667
$ws.Range("N3").Value = -7244.1113

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 4564.1763
$ws.Range("I105").Value = 5007.6665
$ws.Range("K105").Value = 5007.6665
$ws.Range("M105").Value = -3260.6665


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 21563.52
$ws.Range("I58").Value = 23038.195
$ws.Range("K58").Value = 23038.195
$ws.Range("M58").Value = -22835.195

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 20635.18
$ws.Range("I132").Value = 13591.959
$ws.Range("J132").Value = 35640.305
$ws.Range("K132").Value = 40775.877
$ws.Range("L132").Value = 106920.915
$ws.Range("M132").Value = -38245.877
$ws.Range("N132").Value = -111980.915

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3629.327
$ws.Range("I134").Value = 1891.6957
$ws.Range("J134").Value = 16951.166
$ws.Range("K134").Value = 5675.0871
$ws.Range("L134").Value = 50853.49800000001
$ws.Range("M134").Value = -3140.0871
$ws.Range("N134").Value = -55923.49800000001

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 21563.52
$ws.Range("I136").Value = 23038.195
$ws.Range("K136").Value = 69114.58499999999
$ws.Range("M136").Value = -66564.58499999999


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 1092.2128
$ws.Range("I5").Value = 523.4783
$ws.Range("J5").Value = 1637.25
$ws.Range("K5").Value = 1570.4349
$ws.Range("L5").Value = 4911.75
$ws.Range("M5").Value = -1458.4349
$ws.Range("N5").Value = -5135.75

# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 1092.2128
$ws.Range("I135").Value = 523.4783
$ws.Range("J135").Value = 1637.25
$ws.Range("K135").Value = 4711.3047
$ws.Range("L135").Value = 14735.25
$ws.Range("M135").Value = -2176.3047
$ws.Range("N135").Value = -19805.25


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 500103.5
$ws.Range("I2").Value = 687576.3
$ws.Range("K2").Value = 687576.3
$ws.Range("M2").Value = -687463.3

# Row 99 (Leve Item ID 19532)
$ws.Range("H99").Value = 16545.445
$ws.Range("I99").Value = 8418.571
$ws.Range("J99").Value = 44989.5
$ws.Range("K99").Value = 8418.571
$ws.Range("L99").Value = 44989.5
$ws.Range("M99").Value = -6172.571
$ws.Range("N99").Value = -49481.5

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 27453.9
$ws.Range("I132").Value = 109999.5
$ws.Range("K132").Value = 329998.5
$ws.Range("M132").Value = -327468.5


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 766032.7
$ws.Range("I7").Value = 945510.7
$ws.Range("J7").Value = 3251.25
$ws.Range("K7").Value = 945510.7
$ws.Range("L7").Value = 3251.25
$ws.Range("M7").Value = -945398.7
$ws.Range("N7").Value = -3475.25

# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 5490.6665
$ws.Range("I40").Value = 5590.2
$ws.Range("J40").Value = 4993
$ws.Range("K40").Value = 5590.2
$ws.Range("L40").Value = 4993
$ws.Range("M40").Value = -5454.2
$ws.Range("N40").Value = -5265

# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 2913.3635
$ws.Range("I100").Value = 1130.875
$ws.Range("K100").Value = 1130.875
$ws.Range("M100").Value = -589.875

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 766032.7
$ws.Range("I126").Value = 945510.7
$ws.Range("J126").Value = 3251.25
$ws.Range("K126").Value = 2836532.1
$ws.Range("L126").Value = 9753.75
$ws.Range("M126").Value = -2834062.1
$ws.Range("N126").Value = -14693.75

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 6618.6895
$ws.Range("I132").Value = 6488.227
$ws.Range("J132").Value = 7028.7144
$ws.Range("K132").Value = 19464.681
$ws.Range("L132").Value = 21086.1432
$ws.Range("M132").Value = -16934.681
$ws.Range("N132").Value = -26146.1432


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 103 (Leve Item ID 18548)
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 4607.5293
$ws.Range("I126").Value = 4911.484
$ws.Range("J126").Value = 1466.6666
$ws.Range("K126").Value = 14734.452
$ws.Range("L126").Value = 4399.9998
$ws.Range("M126").Value = -12264.452
$ws.Range("N126").Value = -9339.9998

